$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to text format so numeric-looking strings are stored as text,
# matching the source workbook (all cells are inline/shared strings, t="s"/"inlineStr").
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '43.482.22'
$ws.Range("D3").Value = '2.317.68'
$ws.Range("E3").Value = '  +2.26%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '312.31'
$ws.Range("E5").Value = '  +1.74%  '
$ws.Range("D6").Value = '102.96'
$ws.Range("E6").Value = '  +6.30%  '
$ws.Range("D7").Value = '0.538'
$ws.Range("E7").Value = '  +2.33%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("D9").Value = '0.535'
$ws.Range("E9").Value = '  +7.95%  '
$ws.Range("D10").Value = '36.09'
$ws.Range("E10").Value = '  +2.59%  '
$ws.Range("D11").Value = '0.0817'
$ws.Range("E11").Value = '  +3.41%  '
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("D13").Value = '7.06'
$ws.Range("E13").Value = '  +2.61%  '
$ws.Range("D14").Value = '2.675.63'
$ws.Range("E14").Value = '  +2.24%  '
$ws.Range("D15").Value = '15.05'
$ws.Range("E15").Value = '  +2.64%  '
$ws.Range("D16").Value = '2.315.13'
$ws.Range("E16").Value = '  +2.03%  '
$ws.Range("E17").Value = '  +2.55%  '
$ws.Range("D18").Value = '43.385.56'
$ws.Range("E18").Value = '  +3.19%  '
$ws.Range("D19").Value = '12.56'
$ws.Range("E19").Value = '  +1.06%  '
$ws.Range("D20").Value = '0.0₃0924'
$ws.Range("E20").Value = '  +2.02%  '
$ws.Range("D21").Value = '6.18'
$ws.Range("E21").Value = '  +2.85%  '
$ws.Range("D22").Value = '68.46'
$ws.Range("E22").Value = '  +0.28%  '
$ws.Range("D23").Value = '243.18'
$ws.Range("E23").Value = '  +2.13%  '
$ws.Range("E24").Value = '  +6.48%  '
$ws.Range("D25").Value = '2.63'
$ws.Range("E25").Value = '  +2.28%  '
$ws.Range("E26").Value = '  +0.32%  '
$ws.Range("D27").Value = '3.99'
$ws.Range("E27").Value = '  -1.33%  '
$ws.Range("D28").Value = '24.86'
$ws.Range("E28").Value = '  +5.16%  '
$ws.Range("D29").Value = '37.65'
$ws.Range("E29").Value = '  +0.71%  '
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").Value = '9.68'
$ws.Range("E30").Value = '  +2.04%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = '2.12'
$ws.Range("E31").Value = '  +0.50%  '
$ws.Range("D32").Value = '167.58'
$ws.Range("E32").Value = '  +4.41%  '
$ws.Range("D33").Value = '5.35'
$ws.Range("E33").Value = '  +2.40%  '
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("D35").Value = '2.58'
$ws.Range("E35").Value = '  +8.91%  '
$ws.Range("D36").Value = '17.84'
$ws.Range("E36").Value = '  +3.35%  '
$ws.Range("E37").Value = '  +1.18%  '
$ws.Range("E38").Value = '  -2.52%  '
$ws.Range("E39").Value = '  +2.92%  '
$ws.Range("D40").Value = '1.87'
$ws.Range("E40").Value = '  +2.60%  '
$ws.Range("D41").Value = '0.117'
$ws.Range("E41").Value = '  +2.06%  '
$ws.Range("D42").Value = '4.33'
$ws.Range("E42").Value = '  +8.65%  '
$ws.Range("E43").Value = '  +3.39%  '
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("D45").Value = '1.982.00'
$ws.Range("E45").Value = '  +1.55%  '
$ws.Range("E46").Value = '  +3.20%  '
$ws.Range("D47").Value = '3.03'
$ws.Range("E47").Value = '  +5.32%  '
$ws.Range("D48").Value = '9.86'
$ws.Range("E48").Value = '  -1.43%  '
$ws.Range("D49").Value = '56.11'
$ws.Range("E49").Value = '  +5.37%  '
$ws.Range("D50").Value = '2.95'
$ws.Range("E50").Value = '  +15.11%  '
$ws.Range("E51").Value = '  +6.88%  '

# Restore normal style on the range (NumberFormat="@" added a style index;
# resetting to Normal drops it back to the original default style).
$ws.Range("D2:E51").Style = "Normal"

